$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new first column; everything currently in A:D shifts right to B:E
$ws.Columns.Item(1).Insert()

# New headers: B1 is the old "ID" header cell, renamed to "MovieId";
# A1 gets the new "id" header (created after MovieId, matching shared-string order)
$ws.Range("B1").Value = "MovieId"
$ws.Range("A1").Value = "id"

# New GUID "id" values for each data row in column A
$ws.Range("A2").Value = "d861be4a-de63-49ba-94e0-57486b060d90"
$ws.Range("A3").Value = "bd7d27f8-2f3d-4044-8117-e9e71e351339"
$ws.Range("A4").Value = "a7730d44-e048-4879-b571-a8a92a94c1be"
$ws.Range("A5").Value = "de44443e-7c36-4a51-8101-be42d0b572a1"
$ws.Range("A6").Value = "86ae4a99-30aa-42f3-bf6e-0d08e535ff7d"
$ws.Range("A7").Value = "a264bd90-22d5-47b7-aefe-0f6df48de7a3"
$ws.Range("A8").Value = "572eedce-7e46-4d32-915a-f07c529fed2d"

# Drop the old uniform cell style from every used cell (target has none)
$ws.Range("A1:F8").ClearFormats()

# Give the new id column a best-fit width (close to Excel's own autofit result)
$ws.Columns.Item(1).AutoFit()

# Move the active selection the way the author left it
[void]$ws.Range("C9").Select()

Write-Host "done"
